$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 2 (shifts rows 3..72 up by one, removing the
# "add paintings" / "This Cat does not exist" table row).
$ws.Rows("2:2").Delete()

# After the delete, Excel leaves the active cell/selection at A2.
$ws.Range("A2").Select()
